# Weekly update: insert a new observation row (Camote, Vega Central Mapocho
# de Santiago) ahead of the existing row 41, shifting all subsequent rows
# down by one. This mirrors the "insert row" edit recorded in the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 41; rows 41..74 shift down to 42..75.
$ws.Rows("41:41").Insert()

# Populate the newly inserted row 41 with the new observation.
$ws.Range("A41").Value = 9
$ws.Range("B41").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C41").Value = "Metropolitana"
$ws.Range("D41").Value = 44634
$ws.Range("E41").Value = 13
$ws.Range("F41").Value = 100114002
$ws.Range("G41").Value = "Camote"
$ws.Range("H41").Value = "Sin especificar"
$ws.Range("I41").Value = "Primera"
$ws.Range("J41").Value = 1150
$ws.Range("K41").Value = 10000
$ws.Range("L41").Value = 11000
$ws.Range("M41").Value = 10500
$ws.Range("N41").Value = "`$/malla 18 kilos"
$ws.Range("O41").Value = "Perú"
$ws.Range("P41").Value = 583
$ws.Range("Q41").Value = 18
$ws.Range("R41").Value = "Hortaliza"

# Match the date cell's number format to the rest of column D (date/time).
$ws.Range("D41").NumberFormat = $ws.Range("D42").NumberFormat
